# Delete the row for "NNG" / "Nanning, China" (row 264), shifting all
# subsequent rows up by one. This mirrors the author's data regeneration
# where that entry was dropped from the source list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(264).Delete()
